# Updates the Price (D) and Volume(1h) (E) columns of the cryptos sheet
# to the latest scraped values, matching the GitHub Actions commit.
# D-column values that look like plain numbers are written with a
# leading apostrophe so Excel stores them as text (matching the
# original inlineStr cell type) instead of converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.314.79'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.869.19'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''235.10'
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '''0.4694'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').Value = '''0.2863'
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('D10').Value = '''21.79'
$ws.Range('E10').Value = '  -1.80%  '
$ws.Range('D11').Value = '''0.07987'
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('D12').Value = '''96.83'
$ws.Range('E12').Value = '  -1.16%  '
$ws.Range('D13').Value = '1.873.47'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').Value = '''0.6891'
$ws.Range('E14').Value = '  +0.88%  '
$ws.Range('D15').Value = '''5.111'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').Value = '''268.89'
$ws.Range('E16').Value = '  -3.19%  '
$ws.Range('D17').Value = '30.337.82'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '''14.15'
$ws.Range('E18').Value = '  +3.89%  '
$ws.Range('D19').Value = '''0.000007755'
$ws.Range('E19').Value = '  +5.50%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = '2.117.65'
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').Value = '''1.000'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').Value = '''5.257'
$ws.Range('E23').Value = '  -1.99%  '
$ws.Range('D24').Value = '''6.215'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').Value = '''9.392'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('D26').Value = '''167.41'
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('D27').Value = '''18.88'
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('D28').Value = '''1.948'
$ws.Range('E28').Value = '  -0.42%  '
$ws.Range('D29').Value = '''1.364'
$ws.Range('E29').Value = '  -1.35%  '
$ws.Range('D30').Value = '''0.09869'
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('D31').Value = '''4.349'
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('D32').Value = '''1.458'
$ws.Range('E32').Value = '  -1.64%  '
$ws.Range('D33').Value = '''4.058'
$ws.Range('E33').Value = '  -0.46%  '
$ws.Range('D34').Value = '''0.04710'
$ws.Range('E34').Value = '  -0.77%  '
$ws.Range('D35').Value = '''1.136'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('D37').Value = '''2.739'
$ws.Range('E37').Value = '  +1.14%  '
$ws.Range('D38').Value = '''0.01879'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').Value = '''2.820'
$ws.Range('E39').Value = '  +7.16%  '
$ws.Range('D40').Value = '''6.250'
$ws.Range('E40').Value = '  -0.60%  '
$ws.Range('D41').Value = '''72.06'
$ws.Range('E41').Value = '  -4.70%  '
$ws.Range('D42').Value = '''1.958'
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').Value = '''0.4176'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '''0.8424'
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('D45').Value = '''1.000'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('D47').Value = '''7.088'
$ws.Range('E47').Value = '  -1.84%  '
$ws.Range('D48').Value = '''9.168'
$ws.Range('E48').Value = '  -1.32%  '
$ws.Range('D49').Value = '''918.73'
$ws.Range('E49').Value = '  -3.24%  '
$ws.Range('D50').Value = '''34.51'
$ws.Range('E50').Value = '  +0.70%  '
$ws.Range('E51').Value = '  +0.92%  '
